# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme"  (currently used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"      (currently used by the Slide Master)
# The authored change swaps the two theme payloads in place, so the Slide
# Master (i.e. every slide in the deck) picks up the plain "Office Theme"
# color palette that used to belong to the Notes Master, while the
# "Integral" color palette moves the other way.
#
# The font scheme and the fill/line/effect format scheme are byte-for-byte
# identical between the two themes, so the only observable difference after
# the swap is the 12-slot theme color scheme (and the cosmetic theme/
# clrScheme "name" attributes, which PowerPoint does not expose for
# scripted editing). We reproduce the swap by writing the "Office Theme"
# RGB values into the Slide Master's ThemeColorScheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = RGB 0   0   0    # dk1      000000
$tcs.Item(2).RGB  = RGB 255 255 255  # lt1      FFFFFF
$tcs.Item(3).RGB  = RGB 68  84  106  # dk2      44546A
$tcs.Item(4).RGB  = RGB 231 230 230  # lt2      E7E6E6
$tcs.Item(5).RGB  = RGB 91  155 213  # accent1  5B9BD5
$tcs.Item(6).RGB  = RGB 237 125 49   # accent2  ED7D31
$tcs.Item(7).RGB  = RGB 165 165 165  # accent3  A5A5A5
$tcs.Item(8).RGB  = RGB 255 192 0    # accent4  FFC000
$tcs.Item(9).RGB  = RGB 68  114 196  # accent5  4472C4
$tcs.Item(10).RGB = RGB 112 173 71   # accent6  70AD47
$tcs.Item(11).RGB = RGB 5   99  193  # hlink    0563C1
$tcs.Item(12).RGB = RGB 149 79  114  # folHlink 954F72
